$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.619.10"
$ws.Range("E2").Value = "  +4.22%  "
$ws.Range("D3").Value = "1.744.76"
$ws.Range("E3").Value = "  +4.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4794"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2696"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("D10").Value = "1.746.04"
$ws.Range("E10").Value = "  +4.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07108"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6162"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.505"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "26.614.78"
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006895"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("D21").Value = "1.968.63"
$ws.Range("E21").Value = "  +4.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.633"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.858"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.342"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  +2.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.817"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.420"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.023"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.768"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.39%  "
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04570"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.616"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9983"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6355"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9488"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "113.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +18.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.453"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.980"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.04%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01508"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.664"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3910"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.723"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1204"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.73%  "
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.931"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.253"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3452"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.51%  "
